# "puit à colonne permanente"
#
# - rename the last sheet ("Feuil1") to "cuisine" (it already holds the
#   kitchen / Ikea notes, the tab name just never matched)
# - add the new "pcp" (puits à colonne permanente) reference row to the
#   "Géothermie" sheet: a source link in A6 and the term itself in G6
# - leave the view/selection on "Géothermie" (now the active tab), with
#   the "cuisine" sheet's selection parked at C10

$wb = $excel.ActiveWorkbook

# --- rename "Feuil1" -> "cuisine" ------------------------------------------
$wsCuisine = $wb.Worksheets.Item("Feuil1")
$wsCuisine.Name = "cuisine"

# park the selection on the cuisine sheet as it ends up when it's no
# longer the active tab
$wsCuisine.Activate()
$wsCuisine.Range("C10").Select()

# --- add the new geothermy reference row -----------------------------------
$wsGeo = $wb.Worksheets.Item("Géothermie")
$wsGeo.Activate()

$wsGeo.Range("A6").Value = "http://www.ecohabitation.com/actualite/nouvelles/comprendre-geothermie"
$wsGeo.Range("G6").Value = "pcp - puits à colonne permanente"

$wsGeo.Range("G7").Select()
